$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 7 entirely (raw data error: duplicate/incorrect drawdown entry).
# This shifts all rows below it up by one, matching the diff where the
# old row 7 (36612, 15, -0.1114, 36633) disappears and rows 8-13 become 7-12.
$ws.Rows(7).Delete()

# Update the selection to match the post-edit state (C12 selected).
$ws.Range("C12").Select()
